$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8..53 down to 9..54
$ws.Rows.Item(8).EntireRow.Insert()

# Populate the new row 8 with the new weekly record
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(8, 3).Value = "Ñuble"
$ws.Cells.Item(8, 4).Value = 44558
$ws.Cells.Item(8, 5).Value = 16
$ws.Cells.Item(8, 6).Value = 100112031
$ws.Cells.Item(8, 7).Value = "Poroto verde"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 60
$ws.Cells.Item(8, 11).Value = 40000
$ws.Cells.Item(8, 12).Value = 41000
$ws.Cells.Item(8, 13).Value = 40500
$ws.Cells.Item(8, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(8, 15).Value = "Región del Maule"
$ws.Cells.Item(8, 16).Value = 1620
$ws.Cells.Item(8, 17).Value = 25
$ws.Cells.Item(8, 18).Value = "Hortaliza"

# Match the date formatting style used by column D in the rest of the sheet
$ws.Cells.Item(8, 4).NumberFormat = $ws.Cells.Item(9, 4).NumberFormat
